$d = $word.ActiveDocument

# Build characters that don't survive the PS source round-trip reliably.
$eacute = [char]0x00E9
$rsquo  = [char]0x2019

# 1. "s. Sur la culasse" -> "s. Sur la culasse,"
$d.Content.Find.Execute("s. Sur la culasse", $true, $false, $false, $false, $false,
                         $true, 1, $false, "s. Sur la culasse,", 2) | Out-Null

# 2. ". Au devant il porte " -> ". Au devant, il porte "
$d.Content.Find.Execute(". Au devant il porte ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ". Au devant, il porte ", 2) | Out-Null

# 3. "batterie pour avoir bien tost exploicté " -> "batterie, pour avoir bien tost exploicté "
$search3  = "batterie pour avoir bien tost exploict" + $eacute + " "
$replace3 = "batterie, pour avoir bien tost exploict" + $eacute + " "
$d.Content.Find.Execute($search3, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace3, 2) | Out-Null

# 4. Insert a single space right after "<del>d</del>" (only that run's text
#    gains the trailing space; the neighbouring runs are left untouched).
$rng4 = $d.Content
$rng4.Find.Execute("force est de <del>d</del>") | Out-Null
$rng4.Collapse(0)
$rng4.InsertAfter(" ")

# 5. "a balle, de 40" -> "a balle de 40"
$d.Content.Find.Execute("a balle, de 40", $true, $false, $false, $false, $false,
                         $true, 1, $false, "a balle de 40", 2) | Out-Null

# 6. "de canon d'advantaige" -> "de canon dadvantaige" (curly apostrophe removed, runs merged)
$search6  = "de canon d" + $rsquo + "advantaige"
$replace6 = "de canon dadvantaige"
$d.Content.Find.Execute($search6, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace6, 2) | Out-Null

Write-Output "done"
